# Update workbook "Facilidades permanentes 2021 - Diaria" with the latest
# daily data for the Facilidad permanente de depósito (FPD) / liquidez (FPL)
# series (new rows for 15-09-2021 through 30-09-2021).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @("15-09-2021", 3036502, 278768),
    @("16-09-2021", 2187300, 297752),
    @("20-09-2021", 3979500, 225559),
    @("21-09-2021", 5470600, 100000),
    @("22-09-2021", 4445000, 0),
    @("23-09-2021", 4335700, 20000),
    @("24-09-2021", 5276400, 61000),
    @("27-09-2021", 4947550, 20000),
    @("28-09-2021", 5552750, 83000),
    @("29-09-2021", 5517750, 147361),
    @("30-09-2021", 5902400, 149733)
)

$startRow = 179

for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $entry = $newData[$i]
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
}
